$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value but keep it stored as text (matches how this sheet
# keeps numeric-looking figures as text so they don't get reformatted).
function Set-TextValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# Insert a new data row above row 8 (pushes existing rows 8-15 down to 9-16),
# then clone the formatting of row 7 (the first data row) onto the new row 8.
$ws.Rows("8:8").Insert()
$ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))

# Fill in the new row's data - a shortage entry for CIDOPHAGE RETARD.
$ws.Range("A8").Value = 2
Set-TextValue $ws.Range("C8") "CIDOPHAGE RETARD 850MG 60 S.R. TABS"
Set-TextValue $ws.Range("H8") "0:4"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "96.00"
Set-TextValue $ws.Range("P8") "192.0000"
Set-TextValue $ws.Range("Q8") "2:0"

# Renumber the sequence column (م) for the rows that shifted down.
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8

# Update the running total (سعر البيع column) to include the new row.
$ws.Range("P15").Value = 352.55

# Refresh the generated-on timestamp in the footer.
$ws.Range("A16").Value = "Saturday, 20 September, 2025 10:17 AM"
